# Update the cryptos worksheet cells per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.554.95'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').Value = '1.912.25'
$ws.Range('E3').Value = '  +5.35%  '
$ws.Range('E4').Value = '  -0.04%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '315.07'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  +1.47%  '
$cell = $ws.Range('D7')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5227'
$cell.Style = $origStyle
$ws.Range('E7').Value = '  +4.89%  '
$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.3961'
$cell.Style = $origStyle
$ws.Range('E8').Value = '  +1.40%  '
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.09695'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  -1.69%  '
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.152'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  +3.98%  '
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '42.04'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  +2.73%  '
$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.541'
$cell.Style = $origStyle
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('E13').Value = '  +3.08%  '
$ws.Range('D14').Value = '1.917.14'
$ws.Range('E14').Value = '  +5.82%  '
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.549'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  +3.86%  '
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '94.78'
$cell.Style = $origStyle
$ws.Range('E17').Value = '  +2.67%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.00001137'
$cell.Style = $origStyle
$ws.Range('E18').Value = '  -0.28%  '
$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06653'
$cell.Style = $origStyle
$ws.Range('E19').Value = '  +0.31%  '
$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '18.21'
$cell.Style = $origStyle
$ws.Range('E20').Value = '  +5.93%  '
$ws.Range('E21').Value = '  -0.04%  '
$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.336'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  +6.52%  '
$ws.Range('D23').Value = '28.647.05'
$ws.Range('E23').Value = '  +1.86%  '
$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.52'
$cell.Style = $origStyle
$ws.Range('E24').Value = '  +2.50%  '
$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.301'
$cell.Style = $origStyle
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.393'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  -1.06%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.694'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  +12.06%  '
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').Value = '2.130.58'
$ws.Range('E28').Value = '  +5.35%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '21.27'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  +2.82%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '158.68'
$cell.Style = $origStyle
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '129.10'
$cell.Style = $origStyle
$ws.Range('E31').Value = '  +1.89%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.111'
$cell.Style = $origStyle
$ws.Range('E32').Value = '  +7.41%  '
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.1085'
$cell.Style = $origStyle
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.754'
$cell.Style = $origStyle
$ws.Range('E34').Value = '  +3.35%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Range('D35')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.638'
$cell.Style = $origStyle
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.935'
$cell.Style = $origStyle
$ws.Range('E36').Value = '  +11.72%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06773'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.02433'
$cell.Style = $origStyle
$ws.Range('E38').Value = '  +3.81%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.269'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  +7.91%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.2234'
$cell.Style = $origStyle
$ws.Range('E40').Value = '  +4.38%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.82'
$cell.Style = $origStyle
$ws.Range('E41').Value = '  +4.23%  '
$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.092'
$cell.Style = $origStyle
$ws.Range('E42').Value = '  +2.72%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.6469'
$cell.Style = $origStyle
$ws.Range('E43').Value = '  +4.36%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.191'
$cell.Style = $origStyle
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = $origStyle
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '13.67'
$cell.Style = $origStyle
$ws.Range('E46').Value = '  +3.87%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.6101'
$cell.Style = $origStyle
$ws.Range('E47').Value = '  +3.45%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.748'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('B49').Value = 'WEMIXTOKEN'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.286'
$cell.Style = $origStyle
$ws.Range('E49').Value = '  +1.52%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.035'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  +4.92%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '125.50'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  +0.98%  '
